$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121, shifting existing rows 121-175 down to 122-176.
$ws.Range("A121:R121").EntireRow.Insert()

# Populate the newly inserted row 121 with the new record.
$ws.Range("A121").Value = 1
$ws.Range("B121").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C121").Value = "Arica y Parinacota"
$ws.Range("D121").Value = "3/16/2023"
$ws.Range("E121").Value = 15
$ws.Range("F121").Value = 100114001
$ws.Range("G121").Value = "Papa"
$ws.Range("H121").Value = "Asterix"
$ws.Range("I121").Value = "1a (cosecha)"
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 13000
$ws.Range("L121").Value = 14000
$ws.Range("M121").Value = 13550
$ws.Range("N121").Value = "$/saco 25 kilos"
$ws.Range("O121").Value = "Región de Los Lagos"
$ws.Range("P121").Value = 542
$ws.Range("Q121").Value = 25
$ws.Range("R121").Value = "Hortaliza"
